$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")

# Widen column A (now holding a "team object" reference instead of a plain
# name) and re-assert column C at the sheet's default width, same as the
# author's manual column-width pass over A:C.
$ws1.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws1.Columns.Item(3).ColumnWidth = 8.333333333333334

# Make Hoja1 the active/selected sheet (it was Hoja2 before), with the new
# selection sitting on O18:O19.
$ws1.Activate()
$ws1.Range("O18:O19").Select()
